$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("snapshot")

$updates = @{
    "K2"  = "2025-12-16T07:25:01.748837+00:00"
    "K3"  = "2025-12-16T07:25:01.748857+00:00"
    "K4"  = "2025-12-16T07:25:01.748867+00:00"
    "K5"  = "2025-12-16T07:25:04.310989+00:00"
    "K6"  = "2025-12-16T07:25:04.311010+00:00"
    "K7"  = "2025-12-16T07:25:04.311018+00:00"
    "K8"  = "2025-12-16T07:25:06.828853+00:00"
    "K9"  = "2025-12-16T07:25:09.470085+00:00"
    "K10" = "2025-12-16T07:25:11.557565+00:00"
    "K11" = "2025-12-16T07:25:13.969091+00:00"
    "K12" = "2025-12-16T07:25:19.583687+00:00"
    "K13" = "2025-12-16T07:25:19.583715+00:00"
    "K14" = "2025-12-16T07:25:22.208897+00:00"
    "K15" = "2025-12-16T07:25:24.736739+00:00"
    "K16" = "2025-12-16T07:25:26.810942+00:00"
    "K17" = "2025-12-16T07:25:29.308974+00:00"
    "K18" = "2025-12-16T07:25:29.309002+00:00"
    "K19" = "2025-12-16T07:25:29.309020+00:00"
    "K20" = "2025-12-16T07:25:29.309037+00:00"
    "K21" = "2025-12-16T07:25:29.309053+00:00"
    "K22" = "2025-12-16T07:25:31.392117+00:00"
    "K23" = "2025-12-16T07:25:31.392141+00:00"
    "K24" = "2025-12-16T07:25:33.406572+00:00"
    "K25" = "2025-12-16T07:25:33.406596+00:00"
    "K26" = "2025-12-16T07:25:33.406608+00:00"
    "K27" = "2025-12-16T07:25:33.406620+00:00"
    "K28" = "2025-12-16T07:25:33.406632+00:00"
    "K29" = "2025-12-16T07:25:35.991591+00:00"
    "K30" = "2025-12-16T07:25:35.991620+00:00"
    "K31" = "2025-12-16T07:25:35.991638+00:00"
    "K32" = "2025-12-16T07:25:38.620879+00:00"
    "K33" = "2025-12-16T07:25:38.620907+00:00"
    "K34" = "2025-12-16T07:25:38.620924+00:00"
    "K35" = "2025-12-16T07:25:40.692638+00:00"
    "K36" = "2025-12-16T07:25:42.814903+00:00"
    "K37" = "2025-12-16T07:25:42.814919+00:00"
    "K38" = "2025-12-16T07:25:47.414963+00:00"
    "K39" = "2025-12-16T07:25:47.415059+00:00"
    "K40" = "2025-12-16T07:25:49.451756+00:00"
    "K41" = "2025-12-16T07:25:49.451774+00:00"
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
